$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title text (October 2016 -> November 2016) ---
$ws.Range("A1").Value = "Table 4.2. Receipts, Average Cost, and Quality of Fossil Fuels: Electric Utilities, 2006 - November 2016"

# --- Update "Rolling 12 Months Ending in October" -> "...November" (old row 57, before insert) ---
$ws.Range("A57").Value = "Rolling 12 Months Ending in November"

# --- Insert a new row at 53 for the "November" monthly data, pushing subsequent rows down ---
$ws.Rows(53).Insert()

# Copy formatting from row 40 (an existing month-data row) into the new row 53
$ws.Range("A40:M40").Copy()
$ws.Range("A53:M53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 53 with the November monthly data
$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 801020
$ws.Range("C53").Value = 41349
$ws.Range("D53").Value = 2.13
$ws.Range("E53").Value = 41.33
$ws.Range("F53").Value = 1.21
$ws.Range("G53").Value = 116.9
$ws.Range("H53").Value = 6595
$ws.Range("I53").Value = 1074
$ws.Range("J53").Value = 9.79
$ws.Range("K53").Value = 60.1
$ws.Range("L53").Value = 0.48
$ws.Range("M53").Value = 89.6

# --- Update Year-to-Date "2014" row values (now row 55, shifted from 54) ---
$ws.Range("B55").Value = 11003837
$ws.Range("C55").Value = 560287
$ws.Range("D55").Value = 2.37
$ws.Range("E55").Value = 46.58
$ws.Range("F55").Value = 1.21
$ws.Range("G55").Value = 97.4
$ws.Range("H55").Value = 85597
$ws.Range("I55").Value = 14126
$ws.Range("J55").Value = 20.9
$ws.Range("K55").Value = 126.68
$ws.Range("L55").Value = 0.43
$ws.Range("M55").Value = 76.6

# --- Update Year-to-Date "2015" row values (now row 56, shifted from 55) ---
$ws.Range("B56").Value = 10263092
$ws.Range("C56").Value = 528926
$ws.Range("D56").Value = 2.26
$ws.Range("E56").Value = 43.8
$ws.Range("F56").Value = 1.17
$ws.Range("G56").Value = 105.3
$ws.Range("H56").Value = 82638
$ws.Range("I56").Value = 13537
$ws.Range("J56").Value = 11.57
$ws.Range("K56").Value = 70.65
$ws.Range("L56").Value = 0.47
$ws.Range("M56").Value = 77.6

# --- Update Year-to-Date "2016" row values (now row 57, shifted from 56) ---
$ws.Range("B57").Value = 8424882
$ws.Range("C57").Value = 432816
$ws.Range("D57").Value = 2.17
$ws.Range("E57").Value = 42.17
$ws.Range("F57").Value = 1.22
$ws.Range("G57").Value = 95.8
$ws.Range("H57").Value = 67444
$ws.Range("I57").Value = 11011
$ws.Range("J57").Value = 9.03
$ws.Range("K57").Value = 55.29
$ws.Range("L57").Value = 0.46
$ws.Range("M57").Value = 76.8

# --- Update Rolling-12-months "2015" row values (now row 59, shifted from 58) ---
$ws.Range("B59").Value = 11324065
$ws.Range("C59").Value = 583367
$ws.Range("D59").Value = 2.29
$ws.Range("E59").Value = 44.45
$ws.Range("F59").Value = 1.17
$ws.Range("G59").Value = 105.6
$ws.Range("H59").Value = 95399
$ws.Range("I59").Value = 15573
$ws.Range("J59").Value = 11.79
$ws.Range("K59").Value = 72.24
$ws.Range("L59").Value = 0.47
$ws.Range("M59").Value = 83.2

# --- Update Rolling-12-months "2016" row values (now row 60, shifted from 59) ---
$ws.Range("B60").Value = 9250421
$ws.Range("C60").Value = 475597
$ws.Range("D60").Value = 2.17
$ws.Range("E60").Value = 42.21
$ws.Range("F60").Value = 1.21
$ws.Range("G60").Value = 97.1
$ws.Range("H60").Value = 74847
$ws.Range("I60").Value = 12220
$ws.Range("J60").Value = 8.98
$ws.Range("K60").Value = 54.98
$ws.Range("L60").Value = 0.44
$ws.Range("M60").Value = 78.8
